$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.353.10"
$ws.Range("E2").Value = "  -2.53%  "

$ws.Range("D3").Value = "1.942.89"
$ws.Range("E3").Value = "  -2.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.10"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7241"
$ws.Range("E6").Value = "  -8.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3345"
$ws.Range("E8").Value = "  -4.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.83"
$ws.Range("E9").Value = "  +3.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07453"
$ws.Range("E10").Value = "  +6.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8208"
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08138"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").Value = "1.940.78"
$ws.Range("E13").Value = "  -2.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.501"
$ws.Range("E14").Value = "  -1.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.44"
$ws.Range("E15").Value = "  -4.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.94"
$ws.Range("E16").Value = "  -2.93%  "

$ws.Range("D17").Value = "30.373.80"
$ws.Range("E17").Value = "  -2.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008380"
$ws.Range("E18").Value = "  +6.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.59"
$ws.Range("E19").Value = "  -6.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.910"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").Value = "2.196.31"
$ws.Range("E21").Value = "  -2.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.995"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.917"
$ws.Range("E25").Value = "  -1.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.15"
$ws.Range("E26").Value = "  -1.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.424"
$ws.Range("E27").Value = "  +4.45%  "

$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1324"
$ws.Range("E29").Value = "  -11.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.573"
$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.350"
$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.473"
$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.269"
$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05319"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.310"
$ws.Range("E35").Value = "  +6.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7624"
$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01996"
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.851"
$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.47"
$ws.Range("E40").Value = "  +2.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.614"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4585"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.051"
$ws.Range("E43").Value = "  -3.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8452"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.24"
$ws.Range("E46").Value = "  -1.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.864"
$ws.Range("E47").Value = "  -0.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.533"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.12"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4222"
$ws.Range("E50").Value = "  -1.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.517"
$ws.Range("E51").Value = "  -0.95%  "

